$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps ---
# Overview!G2 and de-de!H2 shared the same old value "2016-08-21 06:46:03"
# and both move to the new "generate report" timestamp "2016-08-21 06:46:39".
$wsOverview.Range("G2").Value = "2016-08-21 06:46:39"
$wsDeDe.Range("H2").Value = "2016-08-21 06:46:39"

# zh-cn!H2 had its own value "2016-08-21 06:45:57" -> "2016-08-21 06:46:34"
$wsZhCn.Range("H2").Value = "2016-08-21 06:46:34"

# --- Column widths ---
# Target stored OOXML width is 17.2159881591797 on Overview columns E & F
# and column C on the zh-cn / de-de sheets (was 13.4101845877511).
# This runtime quantizes ColumnWidth to steps of 1/6 of a character, so the
# closest reachable stored width is 17.166666666666668; 16.3 lands there.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
